# fix import dan crud soal kuis reguler
# Cleans up the quiz question import template:
#  - header option labels are lower-cased (opsi_A/B/C/D -> opsi_a/b/c/d)
#  - "jawaban" for the "Pilihan Ganda" graph question now stores the option
#    letter (A) instead of a raw number (4)
#  - the "A. "/"B. "/"C. "/"D. " prefixes are stripped from the answer-option
#    text since the column header already conveys which option it is
#  - the attached image for the population-data question is renamed
#  - a trailing selection is left on H1, matching the re-saved workbook

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header labels -> lower-case opsi_a/b/c/d
$ws.Range("E1").Value = "opsi_a"
$ws.Range("F1").Value = "opsi_b"
$ws.Range("G1").Value = "opsi_c"
$ws.Range("H1").Value = "opsi_d"

# Row 3: "Berapakah rata-rata dari 2, 4, 6?" answer becomes option letter A
$ws.Range("D3").Value = "A"

# Row 4: strip "A. "/"B. "/"C. "/"D. " prefixes from the graph options
$ws.Range("E4").Value = "Grafik A"
$ws.Range("F4").Value = "Grafik B"
$ws.Range("G4").Value = "Grafik C"
$ws.Range("H4").Value = "Grafik D"

# Row 6: rename attached image and strip prefixes from population options
$ws.Range("B6").Value = "statistika.png"
$ws.Range("E6").Value = "200 ribu"
$ws.Range("F6").Value = "300 ribu"
$ws.Range("G6").Value = "400 ribu"
$ws.Range("H6").Value = "500 ribu"

# Match the active selection left behind in the re-saved workbook
$ws.Range("H1").Select()
